$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '255.17'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.95%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.51'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-7.61%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.188'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.21%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05857'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.95%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.712'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.94%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8688'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.35%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9682'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '13.68%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1411'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.03%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07160'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.98%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03181'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.21%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09223'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.39%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001546'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.91%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006079'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006049'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.03%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.498'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.82%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.228'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.57%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.58%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3178'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.76%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.10%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.564'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.28%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04189'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.00%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1399'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.78%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001223'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.07%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004791'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '15.10%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.01%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03816'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.65%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005670'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '58.26%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1103'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2.85%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002300'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-6.50%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009797'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.50%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005372'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.35%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.02%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '21.39%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002127'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-3.67%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.02%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.02%'

Write-Host "Updated symbol list"
